$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.319.51"
$ws.Range("E2").Value = "'  -2.95%  "

$ws.Range("D3").Value = "'2.961.63"
$ws.Range("E3").Value = "'  -2.01%  "

$ws.Range("E4").Value = "'  +0.33%  "

$ws.Range("D5").Value = "'535.90"
$ws.Range("E5").Value = "'  -3.61%  "

$ws.Range("D6").Value = "'148.78"
$ws.Range("E6").Value = "'  -4.69%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.25%  "

$ws.Range("D8").Value = "'0.561"
$ws.Range("E8").Value = "'  +0.80%  "

$ws.Range("D9").Value = "'2.970.59"
$ws.Range("E9").Value = "'  -1.88%  "

$ws.Range("E10").Value = "'  -0.97%  "

$ws.Range("D11").Value = "'6.08"
$ws.Range("E11").Value = "'  -5.43%  "

$ws.Range("D12").Value = "'0.364"
$ws.Range("E12").Value = "'  -0.25%  "

$ws.Range("D13").Value = "'3.483.49"
$ws.Range("E13").Value = "'  -1.71%  "

$ws.Range("E14").Value = "'  +1.54%  "

$ws.Range("D15").Value = "'61.460.85"
$ws.Range("E15").Value = "'  -2.57%  "

$ws.Range("D16").Value = "'23.63"
$ws.Range("E16").Value = "'  -1.70%  "

$ws.Range("D17").Value = "'2.976.65"
$ws.Range("E17").Value = "'  -0.94%  "

$ws.Range("D18").Value = "'0.0000145"
$ws.Range("E18").Value = "'  -3.33%  "

$ws.Range("D19").Value = "'5.11"
$ws.Range("E19").Value = "'  +0.52%  "

$ws.Range("D20").Value = "'11.89"
$ws.Range("E20").Value = "'  -1.16%  "

$ws.Range("D21").Value = "'371.86"
$ws.Range("E21").Value = "'  -6.61%  "

$ws.Range("D22").Value = "'6.66"
$ws.Range("E22").Value = "'  -0.75%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "'  -0.45%  "

$ws.Range("D24").Value = "'65.40"
$ws.Range("E24").Value = "'  +0.77%  "

$ws.Range("D25").Value = "'3.104.20"
$ws.Range("E25").Value = "'  -1.56%  "

$ws.Range("D26").Value = "'0.465"
$ws.Range("E26").Value = "'  -0.15%  "

$ws.Range("D27").Value = "'0.187"
$ws.Range("E27").Value = "'  -0.22%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "'  +0.01%  "

$ws.Range("D29").Value = "'0.0₃0903"
$ws.Range("E29").Value = "'  -7.37%  "

$ws.Range("D30").Value = "'8.08"
$ws.Range("E30").Value = "'  -7.30%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.06%  "

$ws.Range("D32").Value = "'1.70"
$ws.Range("E32").Value = "'  -3.37%  "

$ws.Range("D33").Value = "'20.24"
$ws.Range("E33").Value = "'  -1.18%  "

$ws.Range("D34").Value = "'159.42"
$ws.Range("E34").Value = "'  -1.80%  "

$ws.Range("D35").Value = "'4.51"
$ws.Range("E35").Value = "'  -4.65%  "

$ws.Range("D36").Value = "'5.82"
$ws.Range("E36").Value = "'  -3.53%  "

$ws.Range("D37").Value = "'1.04"
$ws.Range("E37").Value = "'  -6.55%  "

$ws.Range("E38").Value = "'  -4.75%  "

$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "'  -4.60%  "

$ws.Range("B40").Value = "'Maker"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'2.396.96"
$ws.Range("E40").Value = "'  -5.86%  "

$ws.Range("B41").Value = "'OKB"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'37.10"
$ws.Range("E41").Value = "'  -1.31%  "

$ws.Range("D42").Value = "'3.84"
$ws.Range("E42").Value = "'  -1.88%  "

$ws.Range("D43").Value = "'0.664"
$ws.Range("E43").Value = "'  -0.56%  "

$ws.Range("D44").Value = "'21.68"
$ws.Range("E44").Value = "'  -4.91%  "

$ws.Range("D45").Value = "'0.0584"
$ws.Range("E45").Value = "'  -2.66%  "

$ws.Range("D46").Value = "'5.16"
$ws.Range("E46").Value = "'  +1.91%  "

$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "'  +0.21%  "

$ws.Range("D48").Value = "'0.0241"
$ws.Range("E48").Value = "'  -3.22%  "

$ws.Range("D49").Value = "'265.59"
$ws.Range("E49").Value = "'  -2.23%  "

$ws.Range("D50").Value = "'0.0941"
$ws.Range("E50").Value = "'  -0.67%  "

$ws.Range("E51").Value = "'  -0.65%  "
